$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 13
# from 2023-10-08 (45207) to 2023-10-09 (45208), keeping the existing
# date number format / style intact.
$ws.Range("C2:C13").Value = 45208
